# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the second data row
# (8f104349-c4d2-4df8-be52-d8076a42e2d6.md) across the Overview, zh-cn and
# de-de sheets, reflecting a freshly (re)generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-04 12:52:29"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-04 12:52:25"
$zhcn.Range("K3").Value = "2016-09-04 12:52:41"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-04 12:52:29"
$dede.Range("K3").Value = "2016-09-04 12:52:48"
